# Applies the price/volume/coin updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.695.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.201.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.58"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.71"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.73%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.36"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.83"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.101"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.524.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.35"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.195.96"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.552.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.17"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.15"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +15.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0803"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.37"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0328"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.28"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.44"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.52%  "
$ws.Range("B45").Value = "WOONetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.482"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +24.25%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.48"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0980"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +12.17%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.66"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.66%  "
